# Updates the cryptos list (Price column D, Volume(1h) column E) for rows 2-51
# to reflect the latest GitHub Actions scrape, matching the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = "28.707.96"; E = "  +2.39%  " },
    @{ Row = 3; D = "1.872.02"; E = "  +2.25%  " },
    @{ Row = 4; D = $null; E = "  +0.27%  " },
    @{ Row = 5; D = "324.59"; E = "  +0.27%  " },
    @{ Row = 6; D = $null; E = "  +0.05%  " },
    @{ Row = 7; D = $null; E = "  -0.92%  " },
    @{ Row = 8; D = "0.3891"; E = "  +0.74%  " },
    @{ Row = 9; D = "0.07876"; E = "  +0.18%  " },
    @{ Row = 10; D = "0.9757"; E = "  +1.87%  " },
    @{ Row = 11; D = "21.82"; E = "  -0.28%  " },
    @{ Row = 12; D = "1.838.99"; E = "  -0.30%  " },
    @{ Row = 13; D = "7.007"; E = "  +1.56%  " },
    @{ Row = 14; D = "5.704"; E = "  +0.44%  " },
    @{ Row = 15; D = "0.06954"; E = "  +1.51%  " },
    @{ Row = 16; D = "88.35"; E = "  +1.42%  " },
    @{ Row = 17; D = $null; E = "  +0.20%  " },
    @{ Row = 18; D = "0.00001002"; E = "  +1.08%  " },
    @{ Row = 19; D = "16.83"; E = "  +1.46%  " },
    @{ Row = 20; D = "1.003"; E = "  +0.15%  " },
    @{ Row = 21; D = "28.700.49"; E = "  +2.35%  " },
    @{ Row = 22; D = "5.271"; E = "  -0.84%  " },
    @{ Row = 23; D = "11.09"; E = "  +1.06%  " },
    @{ Row = 24; D = "2.105"; E = "  +0.55%  " },
    @{ Row = 25; D = "2.065.90"; E = "  -2.31%  " },
    @{ Row = 26; D = "152.73"; E = "  -0.68%  " },
    @{ Row = 27; D = $null; E = "  +1.16%  " },
    @{ Row = 28; D = "5.869"; E = "  +3.40%  " },
    @{ Row = 29; D = "1.988"; E = "  +1.70%  " },
    @{ Row = 30; D = "119.23"; E = "  +1.43%  " },
    @{ Row = 31; D = "0.09330"; E = "  +0.85%  " },
    @{ Row = 32; D = "0.9197"; E = "  -1.35%  " },
    @{ Row = 33; D = "5.293"; E = "  +0.56%  " },
    @{ Row = 34; D = $null; E = "  +1.22%  " },
    @{ Row = 35; D = "3.321"; E = "  +0.80%  " },
    @{ Row = 36; D = "0.05795"; E = "  -0.61%  " },
    @{ Row = 37; D = "1.154"; E = "  +1.88%  " },
    @{ Row = 38; D = "0.02078"; E = "  -1.86%  " },
    @{ Row = 39; D = "7.679"; E = "  -1.49%  " },
    @{ Row = 40; D = "0.5631"; E = "  +0.80%  " },
    @{ Row = 41; D = "0.1784"; E = "  +1.38%  " },
    @{ Row = 42; D = $null; E = "  -0.82%  " },
    @{ Row = 43; D = "0.07222"; E = "  -0.53%  " },
    @{ Row = 44; D = "11.67"; E = "  +0.87%  " },
    @{ Row = 45; D = "0.5287"; E = "  +0.51%  " },
    @{ Row = 46; D = "2.163"; E = "  +2.06%  " },
    @{ Row = 47; D = "1.126"; E = "  -0.52%  " },
    @{ Row = 48; D = $null; E = "  +0.65%  " },
    @{ Row = 49; D = "112.83"; E = "  +0.46%  " },
    @{ Row = 50; D = "2.407"; E = "  +3.61%  " },
    @{ Row = 51; D = "1.002"; E = "  +0.13%  " }
)

foreach ($u in $updates) {
    $row = $u.Row

    if ($null -ne $u.D) {
        $cellD = $ws.Cells.Item($row, 4)
        $cellD.NumberFormat = "@"
        $cellD.Value2 = $u.D
        $cellD.Style = "Normal"
    }

    $cellE = $ws.Cells.Item($row, 5)
    $cellE.NumberFormat = "@"
    $cellE.Value2 = $u.E
    $cellE.Style = "Normal"
}
